# Regenerate the "K" (strikeouts) column (column G) values in the save-data
# sheet. The column header was previously "Strike#" and is now "K"; the
# underlying per-game strikeout counts were recalculated (std/mean regen)
# and are rewritten here game-by-game.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column G ("K")
$updates = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 0
    6  = 2
    7  = 1
    8  = 1
    9  = 0
    10 = 0
    11 = 3
    12 = 1
    13 = 2
    14 = 0
    15 = 0
    16 = 2
    17 = 0
    18 = 1
    19 = 1
    20 = 0
    21 = 0
    22 = 2
    23 = 1
    25 = 1
    26 = 1
    27 = 1
    28 = 1
    29 = 2
    31 = 2
    32 = 1
    33 = 0
    34 = 1
    35 = 1
    37 = 0
    38 = 0
    39 = 1
    40 = 0
    41 = 0
    42 = 1
    43 = 0
    44 = 0
    45 = 0
    46 = 0
    47 = 1
    48 = 0
    49 = 2
    50 = 1
    51 = 1
    52 = 2
    53 = 1
    54 = 0
    55 = 0
    56 = 3
    57 = 2
    58 = 1
    59 = 0
    60 = 1
    61 = 0
    62 = 2
    63 = 1
    65 = 0
    66 = 2
    67 = 4
    68 = 2
    69 = 2
    70 = 5
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
